# Update "想去人数" (want-to-go count) values in column F across sheets
# to reflect the latest scraped snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 587
$ws.Range("F5").Value = 2573
$ws.Range("F7").Value = 171
$ws.Range("F10").Value = 5306
$ws.Range("F11").Value = 99
$ws.Range("F12").Value = 1464
$ws.Range("F13").Value = 1381
$ws.Range("F14").Value = 594
$ws.Range("F15").Value = 6967
$ws.Range("F17").Value = 48
$ws.Range("F20").Value = 4671
$ws.Range("F24").Value = 1256
$ws.Range("F25").Value = 443
$ws.Range("F26").Value = 1155
$ws.Range("F28").Value = 92
$ws.Range("F29").Value = 70
$ws.Range("F32").Value = 1274
$ws.Range("F34").Value = 238
$ws.Range("F35").Value = 520
$ws.Range("F39").Value = 88
$ws.Range("F42").Value = 1118
$ws.Range("F43").Value = 2405

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 271

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 534
$ws.Range("F8").Value = 1292
$ws.Range("F10").Value = 1740
$ws.Range("F11").Value = 2184
$ws.Range("F13").Value = 518

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 587
$ws.Range("F6").Value = 534
$ws.Range("F7").Value = 2573
$ws.Range("F8").Value = 171
$ws.Range("F9").Value = 1292
$ws.Range("F10").Value = 2184
$ws.Range("F11").Value = 5306
$ws.Range("F15").Value = 99
$ws.Range("F17").Value = 1464
$ws.Range("F18").Value = 1381
$ws.Range("F19").Value = 594
$ws.Range("F20").Value = 6967
$ws.Range("F22").Value = 518
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 4671
$ws.Range("F26").Value = 1256
$ws.Range("F27").Value = 443
$ws.Range("F28").Value = 1155
$ws.Range("F30").Value = 70
$ws.Range("F31").Value = 271
$ws.Range("F37").Value = 238
$ws.Range("F38").Value = 520
$ws.Range("F44").Value = 1118
$ws.Range("F45").Value = 2405
